$wb = $excel.ActiveWorkbook

# "Ready for handoff" -> "Handoff transform failed" everywhere it is used
# (Overview sheet B2/C2, and the per-language sheets' B2 "Status" cell).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handoff transform failed"
$overview.Range("C2").Value = "Handoff transform failed"

$langSheets = @("zh-cn", "de-de")
foreach ($name in $langSheets) {
    $ws = $wb.Worksheets.Item($name)

    # Status column (B) for the handed-off file
    $ws.Range("B2").Value = "Handoff transform failed"

    # Drop the "Latest Handoff File" hyperlink + cell entirely (no more
    # handoff file to link to once the transform failed).
    $ws.Range("C2").Hyperlinks.Delete()
    $ws.Range("C2").Clear()

    # Latest Handoff Datetime reverts to the zero-value sentinel.
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # Latest Handback DateTime sentinel (unchanged text, kept explicit).
    $ws.Range("G2").Value = "0001-01-01 00:00:00"
    $ws.Range("G3").Value = "0001-01-01 00:00:00"
    $ws.Range("D3").Value = "0001-01-01 00:00:00"

    # Handoff Reason: no longer "Include" - this file is ignored.
    $ws.Range("H2").Value = "Ignored"
    $ws.Range("H3").Value = "Ignored"
}
